$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record (dated 2022-06-14, serial 44726) was reported for
# "Femacal de La Calera - Zanahoria". It is inserted as a new row 257,
# which pushes every existing row from 257 downward by one position
# (old row 257 -> new row 258, ..., old row 366 -> new row 367).

# Insert a blank row at position 257, shifting rows 257:366 down to 258:367.
$ws.Rows("257:257").Insert()

# The row that used to be at 257 is now at 258. Use it as a template for
# the brand-new row 257 (same market/region/product/quality/unit/origin
# metadata), then overwrite the measurement columns with the new record's
# values.
$ws.Range("A258:R258").Copy()
$ws.Range("A257:R257").PasteSpecial()

$ws.Range("D257").Value = 44726
$ws.Range("J257").Value = 570
$ws.Range("K257").Value = 6800
$ws.Range("L257").Value = 7000
$ws.Range("M257").Value = 6902
$ws.Range("P257").Value = 345
